{"js": "// Fix three typos in the \"cahier de charge\" document and relocate the\n// stray \"_GoBack\" bookmark left over from the author's last cursor\n// position (Word moves \"_GoBack\" to wherever text was last edited).\nconst doc = context.document;\nconst body = doc.body;\n\n// The document already contains a leftover \"_GoBack\" bookmark further\n// down (an empty paragraph after \"Parcourt du patient\"). Remove it first\n// so we don't end up with two bookmarks sharing the same name once we\n// drop the new one in next to the freshly-edited text.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 1) \"hopital\" -> \"h\u00f4pital\"\nconst hopitalHits = body.search(\"hopital\", { matchCase: true, matchWholeWord: true });\nawait context.sync();\nif (hopitalHits.items.length > 0) {\n  const fixed = hopitalHits.items[0].insertText(\"h\u00f4pital\", \"Replace\");\n  await context.sync();\n  // Word re-drops \"_GoBack\" right after the last text it touched - mirror\n  // that here, immediately following the corrected word.\n  fixed.getRange(\"End\").insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2) \"plusieur\" -> \"plusieurs\" (missing final \"s\")\nconst plusieurHits = body.search(\"plusieur\", { matchCase: true, matchWholeWord: true });\nawait context.sync();\nif (plusieurHits.items.length > 0) {\n  plusieurHits.items[0].getRange(\"End\").insertText(\"s\", \"Replace\");\n  await context.sync();\n}\n\n// 3) \"mogue\" -> \"morgue\" (missing \"r\")\nconst mogueHits = body.search(\"mogue\", { matchCase: true, matchWholeWord: true });\nawait context.sync();\nif (mogueHits.items.length > 0) {\n  mogueHits.items[0].insertText(\"morgue\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Fix three typos in the \"cahier de charge\" document and relocate the\n# stray \"_GoBack\" bookmark left over from the author's last cursor\n# position (Word moves \"_GoBack\" to wherever text was last edited).\n$d = $word.ActiveDocument\n\n# The document already contains a leftover \"_GoBack\" bookmark further\n# down (an empty paragraph after \"Parcourt du patient\"). Remove it first\n# so we don't end up with two bookmarks sharing the same name once we\n# drop the new one in next to the freshly-edited text.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 1) \"hopital\" -> \"h\u00f4pital\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"hopital\", $true, $true)\nif ($found) {\n    $rng.Text = \"h\u00f4pital\"\n    # Word re-drops \"_GoBack\" right after the last text it touched -\n    # mirror that here, immediately following the corrected word.\n    $bmRange = $d.Range($rng.Start, $rng.Start + 7)\n    $bmRange.Collapse(0)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n\n# 2) \"plusieur\" -> \"plusieurs\" (missing final \"s\")\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\"plusieur\", $true, $true)\nif ($found2) {\n    $endRng = $d.Range($rng2.End, $rng2.End)\n    $endRng.InsertAfter(\"s\")\n}\n\n# 3) \"mogue\" -> \"morgue\" (missing \"r\")\n$rng3 = $d.Content\n$found3 = $rng3.Find.Execute(\"mogue\", $true, $true)\nif ($found3) {\n    $rng3.Text = \"morgue\"\n}\n"}
